# Update the Login sheet test data (new generated username/password values)
# and move the active selection, reflecting the refreshed test run described
# in the commit message ("Excel is working now. Updated chrome driver").

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Login")

$ws.Range("B2").Value = "mngr332130"
$ws.Range("C2").Value = "EdAbAda"

# Move / record the active selection on the sheet.
$ws.Activate()
$ws.Range("C6").Select()
